$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Header row
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Surname"
$ws.Range("C1").Value = "Lastname"
$ws.Range("D1").Value = "company"

# Data rows
$ws.Range("A2").Value = "Mr."
$ws.Range("B2").Value = "Lead"
$ws.Range("C2").Value = "Peter"
$ws.Range("D2").Value = "Google"

$ws.Range("A3").Value = "Dr."
$ws.Range("B3").Value = "Friend"
$ws.Range("C3").Value = "Cris"
$ws.Range("D3").Value = "Amazon"

$ws.Range("A4").Value = "Mrs."
$ws.Range("B4").Value = "Lead"
$ws.Range("C4").Value = "Alagu"
$ws.Range("D4").Value = "Ebay"

# Selection as left by the author (cell E6 selected on the contacts sheet)
$ws.Range("E6").Select()
